# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 19 de Marzo de 2020 a las 21:24"

# Swap Israel/Brasil rows (Israel now listed before Brasil) and refresh their figures
$ws.Range("A27").Value = "Israel"
$ws.Range("B27").Value = 677
$ws.Range("C27").Value = 244
$ws.Range("D27").Value = 14
$ws.Range("E27").Value = 663
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 0

$ws.Range("A28").Value = "Brasil"
$ws.Range("B28").Value = 621
$ws.Range("C28").Value = 92
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 612
$ws.Range("F28").Value = 18
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 7

# Update Marruecos (Morocco) figures on row 83
$ws.Range("B83").Value = 63
$ws.Range("C83").Value = 9
$ws.Range("E83").Value = 59
